# Generate Report for Archive
#
# The localization status report is refreshed: every cell whose status was
# "Ready for handoff" moves to "In Translation" (Overview!E2 + Overview!F2,
# the per-language summary cells, and the Status cell on each language
# sheet). Because the new text is shorter than the old text, the "Status"
# column on each sheet is narrower after the refresh.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# Update every cell that currently shows the old status text.
if ($overview.Range("E2").Value() -eq $oldStatus) { $overview.Range("E2").Value = $newStatus }
if ($overview.Range("F2").Value() -eq $oldStatus) { $overview.Range("F2").Value = $newStatus }
if ($zhcn.Range("C2").Value() -eq $oldStatus)      { $zhcn.Range("C2").Value = $newStatus }
if ($dede.Range("C2").Value() -eq $oldStatus)      { $dede.Range("C2").Value = $newStatus }

# The Status column narrows to fit the new (shorter) text.
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
